# planetary_positions.xlsx - alignment fix for Sheet 1 / Sheet 2
# (see commit message: "Fixing of alignment of excel generation rows and
# column Sheet1 and Sheet 2 update.")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: drop the "Name" row, turn the Key/Value header into a single
# merged "Table 1" title row, and fix up the Date value.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet 1")

# The old row 2 ("Name" / "Shreevathsa") goes away; everything below it
# shifts up by one row.
$ws1.Rows.Item(2).Delete()

# Row 1 becomes a single title cell instead of the Key/Value header.
$ws1.Range("B1").Value = $null
$ws1.Range("A1:K1").ClearFormats()
$ws1.Range("A1").Value = "Table 1"
$ws1.Range("A1:K1").Merge()

# Row 2 is now the old "Date" row; refresh its value, forcing text so
# Excel doesn't reinterpret it as a date serial.
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "09/09/1989"
$ws1.Range("B2").ClearFormats()

# ---------------------------------------------------------------------
# Sheet 2: the planetary-position table values were recomputed; update
# each data row (3-14) in place, column by column.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet 2")

# Row 3 - Sun
$ws2.Cells.Item(3, 2).Value = "Leo"
$ws2.Cells.Item(3, 3).Value = "Sun"
$ws2.Cells.Item(3, 4).Value = "Purva Phalguni"
$ws2.Cells.Item(3, 6).Value = 142.4845860305107

# Row 4 - Moon
$ws2.Cells.Item(4, 2).Value = "Scorpio"
$ws2.Cells.Item(4, 3).Value = "Mars"
$ws2.Cells.Item(4, 4).Value = "Jyeshtha"
$ws2.Cells.Item(4, 5).Value = "Mercury"
$ws2.Cells.Item(4, 6).Value = 238.080245280088
$ws2.Cells.Item(4, 10).Value = 5

# Row 5 - Mercury
$ws2.Cells.Item(5, 2).Value = "Virgo"
$ws2.Cells.Item(5, 3).Value = "Mercury"
$ws2.Cells.Item(5, 4).Value = "Hasta"
$ws2.Cells.Item(5, 5).Value = "Moon"
$ws2.Cells.Item(5, 6).Value = 166.4089687091878
$ws2.Cells.Item(5, 10).Value = 3

# Row 6 - Venus
$ws2.Cells.Item(6, 2).Value = "Libra"
$ws2.Cells.Item(6, 3).Value = "Venus"
$ws2.Cells.Item(6, 4).Value = "Chitra"
$ws2.Cells.Item(6, 5).Value = "Mars"
$ws2.Cells.Item(6, 6).Value = 182.0572724410369

# Row 7 - Mars
$ws2.Cells.Item(7, 2).Value = "Leo"
$ws2.Cells.Item(7, 3).Value = "Sun"
$ws2.Cells.Item(7, 4).Value = "Uttara Phalguni"
$ws2.Cells.Item(7, 5).Value = "Sun"
$ws2.Cells.Item(7, 6).Value = 149.3837157754606
$ws2.Cells.Item(7, 7).Value = "Direct"
$ws2.Cells.Item(7, 8).Value = "Combust"
$ws2.Cells.Item(7, 10).Value = 2

# Row 8 - Jupiter
$ws2.Cells.Item(8, 2).Value = "Gemini"
$ws2.Cells.Item(8, 3).Value = "Mercury"
$ws2.Cells.Item(8, 4).Value = "Ardra"
$ws2.Cells.Item(8, 5).Value = "Rahu"
$ws2.Cells.Item(8, 6).Value = 73.33655250974485
$ws2.Cells.Item(8, 10).Value = 12

# Row 9 - Saturn
$ws2.Cells.Item(9, 2).Value = "Sagittarius"
$ws2.Cells.Item(9, 3).Value = "Jupiter"
$ws2.Cells.Item(9, 4).Value = "Purva Ashadha"
$ws2.Cells.Item(9, 6).Value = 253.590182561376
$ws2.Cells.Item(9, 7).Value = "Retro"
$ws2.Cells.Item(9, 8).Value = "No"
$ws2.Cells.Item(9, 10).Value = 6

# Row 10 - Uranus
$ws2.Cells.Item(10, 2).Value = "Sagittarius"
$ws2.Cells.Item(10, 3).Value = "Jupiter"
$ws2.Cells.Item(10, 4).Value = "Moola"
$ws2.Cells.Item(10, 5).Value = "Ketu"
$ws2.Cells.Item(10, 6).Value = 247.622252771897
$ws2.Cells.Item(10, 7).Value = "Retro"
$ws2.Cells.Item(10, 10).Value = 6

# Row 11 - Neptune
$ws2.Cells.Item(11, 2).Value = "Sagittarius"
$ws2.Cells.Item(11, 3).Value = "Jupiter"
$ws2.Cells.Item(11, 4).Value = "Purva Ashadha"
$ws2.Cells.Item(11, 5).Value = "Venus"
$ws2.Cells.Item(11, 6).Value = 255.9345667638737
$ws2.Cells.Item(11, 10).Value = 6

# Row 12 - Pluto
$ws2.Cells.Item(12, 2).Value = "Libra"
$ws2.Cells.Item(12, 3).Value = "Venus"
$ws2.Cells.Item(12, 4).Value = "Swati"
$ws2.Cells.Item(12, 5).Value = "Rahu"
$ws2.Cells.Item(12, 6).Value = 199.2862020569684
$ws2.Cells.Item(12, 7).Value = "Direct"
$ws2.Cells.Item(12, 10).Value = 4

# Row 13 - Rahu
$ws2.Cells.Item(13, 2).Value = "Aquarius"
$ws2.Cells.Item(13, 3).Value = "Saturn"
$ws2.Cells.Item(13, 4).Value = "Dhanishta"
$ws2.Cells.Item(13, 5).Value = "Mars"
$ws2.Cells.Item(13, 6).Value = 300.7872291398033
$ws2.Cells.Item(13, 10).Value = 8

# Row 14 - Ketu
$ws2.Cells.Item(14, 2).Value = "Leo"
$ws2.Cells.Item(14, 3).Value = "Sun"
$ws2.Cells.Item(14, 4).Value = "Magha"
$ws2.Cells.Item(14, 5).Value = "Ketu"
$ws2.Cells.Item(14, 6).Value = 121.9004845901756
$ws2.Cells.Item(14, 7).Value = "Direct"
$ws2.Cells.Item(14, 10).Value = 2
